$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# New shared-string texts (appended in this exact order so they end
# up at shared-string indices 22-34, matching the target workbook).
# -----------------------------------------------------------------

$a8 = 'Criar tabela SegmentoClassificacao'
$b8 = 'Preciso gerar um script para criar uma tabela no SQL Server onde deverá constar a primary key identity ID, a sigla e o descritivo.'

$a9 = 'Criar tabela Segmento'
$b9 = 'Preciso gerar um script para criar uma tabela no SQL Server com o nome "Setor Econômico" onde deverá constar a primary key identity ID e o descritivo.'

$a10 = 'Criar tabela Subsetor'
$b10 = 'Preciso gerar um script para criar uma tabela no SQL Server com o nome "Subsetor" onde deverá constar a primary key identity ID e o descritivo.'

$a11 = 'Criar tabela SetorEconomico'
$b11 = 'Preciso gerar um script para criar uma tabela no SQL Server com o nome "Segmento" onde deverá constar a primary key identity ID e o descritivo.'

$a12 = 'Criar tabela Empresa'
$b12 = @'
Preciso gerar um script para criar uma tabela no SQL Server com o nome "Empresa" onde deverá constar a primary key identity ID, Nome, código (sigla) fazendo chave estrangeira para as seguintes 
tabelas dbo.SegmentoClassificacao (NÃO obrigatória), SetorEconomico (obrigatória), Subsetor (obrigatória) e Segmento (obrigatória).
'@

$c8 = @'
Aqui está um exemplo de script SQL para criar uma tabela no SQL Server com as colunas que você mencionou:
SQL
CREATE TABLE [dbo].[Tabela_Siglas] (
    [ID] INT IDENTITY(1,1) PRIMARY KEY,
    [Sigla] VARCHAR(10) NOT NULL,
    [Descritivo] VARCHAR(100) NOT NULL
);
Esse script cria uma tabela chamada Tabela_Siglas com as seguintes colunas:
ID: uma coluna de tipo INT que serve como chave primária e é auto-incrementada usando a propriedade IDENTITY(1,1).
Sigla: uma coluna de tipo VARCHAR(10) que armazena a sigla.
Descritivo: uma coluna de tipo VARCHAR(100) que armazena o descritivo.
Você pode executar esse script no SQL Server Management Studio ou em qualquer outro cliente SQL para criar a tabela.
Lembre-se de que você pode ajustar os tipos de dados e as restrições de acordo com as necessidades específicas do seu aplicativo.
'@

$c9 = @'
Aqui está um exemplo de script SQL para criar uma tabela no SQL Server com as colunas que você mencionou:
SQL
CREATE TABLE [dbo].[Setor_Economico] (
    [ID] INT IDENTITY(1,1) PRIMARY KEY,
    [Descritivo] VARCHAR(100) NOT NULL
);
Esse script cria uma tabela chamada Setor_Economico com as seguintes colunas:
ID: uma coluna de tipo INT que serve como chave primária e é auto-incrementada usando a propriedade IDENTITY(1,1).
Descritivo: uma coluna de tipo VARCHAR(100) que armazena o descritivo do setor econômico.
Você pode executar esse script no SQL Server Management Studio ou em qualquer outro cliente SQL para criar a tabela.
Lembre-se de que você pode ajustar os tipos de dados e as restrições de acordo com as necessidades específicas do seu aplicativo.
'@

$c10 = @'
Aqui está um exemplo de script SQL para criar uma tabela no SQL Server com as colunas que você mencionou:
SQL
CREATE TABLE [dbo].[Subsetor] (
    [ID] INT IDENTITY(1,1) PRIMARY KEY,
    [Descritivo] VARCHAR(100) NOT NULL
);
Esse script cria uma tabela chamada Subsetor com as seguintes colunas:
ID: uma coluna de tipo INT que serve como chave primária e é auto-incrementada usando a propriedade IDENTITY(1,1).
Descritivo: uma coluna de tipo VARCHAR(100) que armazena o descritivo do subsetor.
Você pode executar esse script no SQL Server Management Studio ou em qualquer outro cliente SQL para criar a tabela.
Lembre-se de que você pode ajustar os tipos de dados e as restrições de acordo com as necessidades específicas do seu aplicativo.
'@

# Strip the single trailing newline that the here-strings add after the
# final line (PowerShell here-strings always end with a line break before
# the closing `'@`), and normalize all line endings to the Unix style
# used by the workbook.
function Normalize([string]$t) {
    $t = $t -replace "`r`n", "`n"
    if ($t.EndsWith("`n")) { $t = $t.Substring(0, $t.Length - 1) }
    return $t
}
$b12 = Normalize $b12
$c8 = Normalize $c8
$c9 = Normalize $c9
$c10 = Normalize $c10

# -----------------------------------------------------------------
# Write column A/B for every new row first, then column C, matching
# the order in which the author originally typed the sheet (titles
# first, long SQL answers last) so new shared-string indices land in
# the same order as the target workbook (22-34).
# -----------------------------------------------------------------
$ws.Range("A8").Value2 = $a8
$ws.Range("B8").Value2 = $b8

$ws.Range("A9").Value2 = $a9
$ws.Range("B9").Value2 = $b9

$ws.Range("A10").Value2 = $a10
$ws.Range("B10").Value2 = $b10

$ws.Range("A11").Value2 = $a11
$ws.Range("B11").Value2 = $b11

$ws.Range("A12").Value2 = $a12
$ws.Range("B12").Value2 = $b12
$ws.Range("B12").WrapText = $true

$ws.Range("C8").Value2 = $c8
$ws.Range("C8").WrapText = $true

$ws.Range("C9").Value2 = $c9
$ws.Range("C9").WrapText = $true

$ws.Range("C10").Value2 = $c10
$ws.Range("C10").WrapText = $true

# -----------------------------------------------------------------
# Row 8 time columns + formula
# -----------------------------------------------------------------
$ws.Range("E8").Value2 = 0.96736111111111112
$ws.Range("E8").NumberFormat = "h:mm"
$ws.Range("F8").Value2 = 0.96944444444444444
$ws.Range("F8").NumberFormat = "h:mm"
$ws.Range("G8").Formula = "=F8-E8"
$ws.Range("G8").NumberFormat = "h:mm"

# -----------------------------------------------------------------
# Row 9 time column (no F/G - task still open)
# -----------------------------------------------------------------
$ws.Range("E9").Value2 = 0.96944444444444444
$ws.Range("E9").NumberFormat = "h:mm"

# -----------------------------------------------------------------
# Row heights
# -----------------------------------------------------------------
$ws.Rows.Item(8).RowHeight = 360
$ws.Rows.Item(9).RowHeight = 315
$ws.Rows.Item(10).RowHeight = 315
$ws.Rows.Item(12).RowHeight = 120

# -----------------------------------------------------------------
# Sheet view: scroll so row 10 is at the top, and select B1:B12
# (mirrors the author's on-screen state after finishing data entry).
# -----------------------------------------------------------------
$excel.Goto($ws.Range("B1:B12"))
$excel.ActiveWindow.ScrollRow = 10
